# live_trading_results.xlsx — record Trade #5 (closed 2026-02-17 13:07:57)
# and roll the new totals up into the Summary / Strategy Status sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.52   # Current Capital
$summary.Range("B4").Value = -0.48     # Total P&L $
$summary.Range("B5").Value = -1.92     # Total P&L %
$summary.Range("B6").Value = 5         # Total Trades
$summary.Range("B8").Value = 4         # Losing Trades
$summary.Range("B9").Value = 20        # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet — MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.52
$status.Range("D4").Value = 5
$status.Range("E4").Value = -0.48
$status.Range("F4").Value = -0.48
$status.Range("G4").Value = 20

# ---------------------------------------------------------------------
# New trade row (#5) appended to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(6, 1).Value = 5                 # Trade #

    # "2026-02-17" reads as a date literal, so Excel would otherwise
    # auto-convert it to a date serial the way it never did for the
    # existing rows (which store it as plain text). Force text, write
    # it, then drop the now-unneeded format so the cell matches its
    # plain, style-less siblings.
    $ws.Cells.Item(6, 2).NumberFormat = "@"
    $ws.Cells.Item(6, 2).Value = "2026-02-17"      # Date
    $ws.Cells.Item(6, 2).ClearFormats()

    $ws.Cells.Item(6, 3).Value = "13:07:50"        # Time
    $ws.Cells.Item(6, 4).Value = "MarketMaking"    # Strategy
    $ws.Cells.Item(6, 5).Value = "UP"              # Side
    $ws.Cells.Item(6, 6).Value = 0.19              # Entry Price
    $ws.Cells.Item(6, 7).Value = 0.14              # Exit Price
    $ws.Cells.Item(6, 8).Value = "CLOSED"          # Status
    $ws.Cells.Item(6, 9).Value = -26.3158          # P&L %
    $ws.Cells.Item(6, 10).Value = -0.05            # P&L $
    $ws.Cells.Item(6, 11).Value = 99.52            # Capital After
    $ws.Cells.Item(6, 12).Value = 0                # Entry Slippage (bps)
    $ws.Cells.Item(6, 13).Value = 0                # Exit Slippage (bps)
    $ws.Cells.Item(6, 14).Value = 0.6              # Confidence
    $ws.Cells.Item(6, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
    $ws.Cells.Item(6, 16).Value = "early_exit"     # Exit Reason
    $ws.Cells.Item(6, 17).Value = 0.11             # Duration (min)
}
